$d = $word.ActiveDocument

# Update the date paragraph
$null = $d.Content.Find.Execute("2025-01-23 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-01-24 Friday", 2)

# Update each table cell value positionally (row-major), since several
# old values are duplicated across cells and must not be replaced via a
# global find/replace.
$t = $d.Tables.Item(1)
$values = @(
    "27+41=68",
    "98-66=32",
    "85-67=18",
    "66-13=53",
    "55-5=50",
    "33-29=4",
    "32+27=59",
    "82-67=15",
    "44-34=10",
    "79-63=16",
    "50-1=49",
    "62-37=25",
    "82-67=15",
    "45+26=71",
    "96-57=39",
    "36-5=31",
    "50+39=89",
    "98-55=43",
    "93+1=94",
    "73-2=71",
    "57+17=74",
    "90-22=68",
    "42+27=69",
    "12+44=56",
    "91-42=49",
    "10+9=19",
    "14+10=24",
    "37-20=17",
    "54-30=24",
    "99-61=38",
    "2+62=64",
    "60-34=26",
    "61+13=74",
    "8+11=19",
    "51-22=29",
    "25+37=62",
    "9+42=51",
    "65-27=38",
    "97-75=22",
    "21+21=42",
    "3+34=37",
    "24+58=82",
    "78+12=90",
    "93-78=15",
    "85-69=16",
    "7+61=68",
    "44-5=39",
    "92-52=40",
    "40-28=12",
    "32+9=41",
    "71-9=62",
    "56+41=97",
    "60-9=51",
    "62-38=24",
    "0+12=12",
    "52-8=44",
    "51+2=53",
    "38+60=98",
    "54-49=5",
    "27+51=78",
    "12-0=12",
    "50-16=34",
    "10+33=43",
    "94-85=9",
    "42+18=60",
    "84-38=46",
    "68-55=13",
    "26+22=48",
    "9+53=62",
    "23-13=10",
    "95-57=38",
    "54-7=47",
    "28+39=67",
    "20+46=66",
    "76-61=15",
    "42+43=85",
    "9+79=88",
    "45+44=89",
    "76-1=75",
    "82-20=62",
    "27+19=46",
    "81-4=77",
    "35-22=13",
    "63+27=90",
    "11+38=49",
    "4+59=63",
    "69+24=93",
    "98-42=56",
    "44+30=74",
    "78-51=27",
    "37-27=10",
    "90-44=46",
    "14+84=98",
    "82-17=65",
    "9+41=50",
    "63-33=30",
    "5-0=5",
    "29+35=64",
    "40+26=66",
    "69+0=69"
)

if ($t.Rows.Count * $t.Columns.Count -ne $values.Length) {
    throw "Unexpected table size: $($t.Rows.Count) rows x $($t.Columns.Count) cols, expected $($values.Length) values"
}

$idx = 0
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    for ($c = 1; $c -le $t.Columns.Count; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $values[$idx]
        $idx = $idx + 1
    }
}

Write-Host "Updated $idx cells"